# #1626 Extend restore clipboard for actions that return a value (#1627)
# Fix convert-to-picture test: the shape created by pasting the copied
# picture was still named "pictocopy" (leftover from the copy source)
# and the caption below it was split across two runs. Rename the
# pasted shape and collapse the caption back into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- Rename the pasted picture (id=6) from "pictocopy" to "copied" ---
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 6 -and $shp.Name -eq "pictocopy") {
        $shp.Name = "copied"
    }
}

# --- Merge the caption's two runs into a single run ---
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "text 3") {
        $tr = $shp.TextFrame.TextRange
        $para = $tr.Paragraphs(1)
        # Force the run-merge: setting the same concatenated text as a
        # no-op change leaves the original run split untouched, so flip
        # through a distinct placeholder value first.
        $para.Text = "__tmp__"
        $para2 = $tr.Paragraphs(1)
        $para2.Text = "Expected Output (pasted object is the right shape )"
    }
}
